$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "level 2 6RES source plate"
$ws.Range("C3").Value = "6RES_AQ_BP"
$ws.Range("E3").Value = "384-Well Level 2 MoClo output plate"
$ws.Range("F3").Value = "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)"
$ws.Range("G3").Value = "A2"
$ws.Range("H3").Value = 2875
$ws.Range("I3").Value = "Deionised water"
